$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Seed new rows 72-77 by copying the formatting (styles) of the last existing
# data row (71) so column A keeps its bold/bordered/centered style (s=1) and
# column E keeps its datetime number format (s=2), matching the rest of the table.
$ws.Range("A71:V71").Copy()
$ws.Range("A72:V77").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Row 72: match index 71
$ws.Range("A72").Value = 71
$ws.Range("B72").Value = "poland"
$ws.Range("C72").Value = "ekstraklasa"
$ws.Range("D72").Value = "2023-2024"
$ws.Range("E72").Value = 45192.625
$ws.Range("F72").Value = "Slask Wroclaw"
$ws.Range("G72").Value = 1
$ws.Range("H72").Value = "Piast Gliwice"
$ws.Range("I72").Value = 0
$ws.Range("J72").Value = 3.16
$ws.Range("K72").Value = "16/09/2023 19:13"
$ws.Range("L72").Value = 4.04
$ws.Range("M72").Value = "23/09/2023 14:54"
$ws.Range("N72").Value = 3.04
$ws.Range("O72").Value = "16/09/2023 19:13"
$ws.Range("P72").Value = 3.12
$ws.Range("Q72").Value = "23/09/2023 14:54"
$ws.Range("R72").Value = 2.51
$ws.Range("S72").Value = "16/09/2023 19:13"
$ws.Range("T72").Value = 2.15
$ws.Range("U72").Value = "23/09/2023 14:54"
$ws.Range("V72").Value = "https://www.betexplorer.com/football/poland/ekstraklasa/slask-wroclaw-piast-gliwice/UXNn1PGk/"

# Row 73: match index 72
$ws.Range("A73").Value = 72
$ws.Range("B73").Value = "poland"
$ws.Range("C73").Value = "ekstraklasa"
$ws.Range("D73").Value = "2023-2024"
$ws.Range("E73").Value = 45192.72916666666
$ws.Range("F73").Value = "Cracovia"
$ws.Range("G73").Value = 1
$ws.Range("H73").Value = "Pogon Szczecin"
$ws.Range("I73").Value = 5
$ws.Range("J73").Value = 2.78
$ws.Range("K73").Value = "17/09/2023 14:13"
$ws.Range("L73").Value = 2.76
$ws.Range("M73").Value = "23/09/2023 16:04"
$ws.Range("N73").Value = 3.35
$ws.Range("O73").Value = "17/09/2023 14:13"
$ws.Range("P73").Value = 3.4
$ws.Range("Q73").Value = "23/09/2023 14:45"
$ws.Range("R73").Value = 2.49
$ws.Range("S73").Value = "17/09/2023 14:13"
$ws.Range("T73").Value = 2.66
$ws.Range("U73").Value = "23/09/2023 16:04"
$ws.Range("V73").Value = "https://www.betexplorer.com/football/poland/ekstraklasa/cracovia-pogon-szczecin/GUcDnQ1d/"

# Row 74: match index 73
$ws.Range("A74").Value = 73
$ws.Range("B74").Value = "poland"
$ws.Range("C74").Value = "ekstraklasa"
$ws.Range("D74").Value = "2023-2024"
$ws.Range("E74").Value = 45192.83333333334
$ws.Range("F74").Value = "Lech Poznan"
$ws.Range("G74").Value = 2
$ws.Range("H74").Value = "Stal Mielec"
$ws.Range("I74").Value = 1
$ws.Range("J74").Value = 1.29
$ws.Range("K74").Value = "18/09/2023 18:12"
$ws.Range("L74").Value = 1.35
$ws.Range("M74").Value = "23/09/2023 19:57"
$ws.Range("N74").Value = 5.65
$ws.Range("O74").Value = "18/09/2023 18:12"
$ws.Range("P74").Value = 5.32
$ws.Range("Q74").Value = "23/09/2023 19:57"
$ws.Range("R74").Value = 10.36
$ws.Range("S74").Value = "18/09/2023 18:12"
$ws.Range("T74").Value = 8.99
$ws.Range("U74").Value = "23/09/2023 19:57"
$ws.Range("V74").Value = "https://www.betexplorer.com/football/poland/ekstraklasa/lech-poznan-stal-mielec/Qau49m99/"

# Row 75: match index 74
$ws.Range("A75").Value = 74
$ws.Range("B75").Value = "poland"
$ws.Range("C75").Value = "ekstraklasa"
$ws.Range("D75").Value = "2023-2024"
$ws.Range("E75").Value = 45193.52083333334
$ws.Range("F75").Value = "Zaglebie"
$ws.Range("G75").Value = 1
$ws.Range("H75").Value = "Warta Poznan"
$ws.Range("I75").Value = 0
$ws.Range("J75").Value = 2.16
$ws.Range("K75").Value = "18/09/2023 18:12"
$ws.Range("L75").Value = 2.05
$ws.Range("M75").Value = "24/09/2023 12:21"
$ws.Range("N75").Value = 3.32
$ws.Range("O75").Value = "18/09/2023 18:12"
$ws.Range("P75").Value = 3.4
$ws.Range("Q75").Value = "24/09/2023 12:21"
$ws.Range("R75").Value = 3.58
$ws.Range("S75").Value = "18/09/2023 18:12"
$ws.Range("T75").Value = 3.98
$ws.Range("U75").Value = "24/09/2023 12:21"
$ws.Range("V75").Value = "https://www.betexplorer.com/football/poland/ekstraklasa/zaglebie-warta-poznan/hnsdB9vc/"

# Row 76: match index 75
$ws.Range("A76").Value = 75
$ws.Range("B76").Value = "poland"
$ws.Range("C76").Value = "ekstraklasa"
$ws.Range("D76").Value = "2023-2024"
$ws.Range("E76").Value = 45193.625
$ws.Range("F76").Value = "Ruch Chorzow"
$ws.Range("G76").Value = 3
$ws.Range("H76").Value = "Rakow"
$ws.Range("I76").Value = 5
$ws.Range("J76").Value = 4.78
$ws.Range("K76").Value = "19/09/2023 13:42"
$ws.Range("L76").Value = 4.08
$ws.Range("M76").Value = "24/09/2023 14:42"
$ws.Range("N76").Value = 3.81
$ws.Range("O76").Value = "19/09/2023 13:42"
$ws.Range("P76").Value = 3.63
$ws.Range("Q76").Value = "24/09/2023 14:42"
$ws.Range("R76").Value = 1.74
$ws.Range("S76").Value = "19/09/2023 13:42"
$ws.Range("T76").Value = 1.94
$ws.Range("U76").Value = "24/09/2023 14:42"
$ws.Range("V76").Value = "https://www.betexplorer.com/football/poland/ekstraklasa/ruch-chorzow-rakow-czestochowa/bkXG6ogS/"

# Row 77: match index 76
$ws.Range("A77").Value = 76
$ws.Range("B77").Value = "poland"
$ws.Range("C77").Value = "ekstraklasa"
$ws.Range("D77").Value = "2023-2024"
$ws.Range("E77").Value = 45193.72916666666
$ws.Range("F77").Value = "Legia"
$ws.Range("G77").Value = 2
$ws.Range("H77").Value = "Gornik Zabrze"
$ws.Range("I77").Value = 1
$ws.Range("J77").Value = 1.5
$ws.Range("K77").Value = "17/09/2023 16:43"
$ws.Range("L77").Value = 1.52
$ws.Range("M77").Value = "24/09/2023 17:26"
$ws.Range("N77").Value = 4.47
$ws.Range("O77").Value = "17/09/2023 16:43"
$ws.Range("P77").Value = 4.41
$ws.Range("Q77").Value = "24/09/2023 17:18"
$ws.Range("R77").Value = 6.45
$ws.Range("S77").Value = "17/09/2023 16:43"
$ws.Range("T77").Value = 6.37
$ws.Range("U77").Value = "24/09/2023 17:18"
$ws.Range("V77").Value = "https://www.betexplorer.com/football/poland/ekstraklasa/legia-gornik-zabrze/xzq887OF/"

$ws.Range("A1").Select()
